# repull data, push all data, mean calculation
# Update column F (dSF) values for rows 2-33 per freshly re-pulled source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 5
    3  = 3
    4  = 1
    5  = -2
    6  = 3
    7  = -5
    8  = 2
    10 = 2
    11 = 5
    12 = 4
    13 = -2
    14 = -5
    15 = 1
    16 = 2
    17 = -1
    18 = 3
    20 = 2
    21 = 2
    22 = 1
    23 = -2
    24 = -3
    25 = -2
    26 = -6
    27 = 1
    28 = -2
    29 = 4
    30 = 5
    31 = 1
    32 = -1
    33 = 9
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
